# Application is now targeting .NET 4.0 Client Profile.
# Update the "DatePicker controls look strange" bug row (row 5):
#  - Status: WONTFIX -> FIXED
#  - Description: append Note1/Note2 about the .NET 4.0 fix
#  - Last modified date bumped to 2010-03-18 (serial 40255)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C5").Value = "FIXED"
$ws.Range("E5").Value = "After updating to the latest WPF Toolkit, the DatePicker controls all look disabled, even when they shouldn't be.  Note1: It appears that this is a bug with WPF Toolkit, so I will just wait until they fix it.  Note2: I have fixed this issue by switching to the DataGrid and DatePicker controls in .net 4.0"
$ws.Range("G5").Value = 40255

# The longer description needs a taller row to keep displaying fully wrapped.
$ws.Rows.Item(5).RowHeight = 60

# Scroll the view a bit to the right and move the selection, as the author
# left the sheet before saving.
$ws.Application.Goto($ws.Range("E6"))
$ws.Range("E6").Select()
